$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# --- Numeric / reused-string cells first (do not allocate new shared strings) ---
$ws.Range("A4").Value = 50000002
$ws.Range("D4").Value = "Upgrade"
$ws.Range("F4").Value = 1234
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 50
$ws.Range("J4").NumberFormat = "#,##0"
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 250
$ws.Range("O4").Value = 10
$ws.Range("P4").Value = 20
$ws.Range("Q4").Value = 30
$ws.Range("R4").Value = 40
$ws.Range("S4").Value = 50

$ws.Range("A5").Value = 50000003
$ws.Range("D5").Value = "Upgrade"
$ws.Range("F5").Value = 1234
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 50
$ws.Range("J5").NumberFormat = "#,##0"
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 250
$ws.Range("O5").Value = 30
$ws.Range("P5").Value = 60
$ws.Range("Q5").Value = 90
$ws.Range("R5").Value = 120
$ws.Range("S5").Value = 150

# --- New text cells, written in the exact order the new labels were authored ---
$ws.Range("C4").Value = "체력 #N 증가"
$ws.Range("B5").Value = "시작 골드 증가"
$ws.Range("C5").Value = "골드 #N 증가"
$ws.Range("E5").Value = "StartGold"
$ws.Range("E4").Value = "MaxHP"
$ws.Range("B4").Value = "최대 체력 증가"

# Update selection to E5 to match the saved cursor position
$ws.Range("E5").Select()
